$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell C4: numeric 2 -> text "Stand" (written first so it becomes shared-string index 0)
$ws.Range("C4").Value = "Stand"

# Cell A3: numeric 3 -> text "Text" (written second so it becomes shared-string index 1)
$ws.Range("A3").Value = "Text"

# Update the active selection shown in the sheet view (J5 -> J9)
$ws.Range("J9").Select() | Out-Null
